# Applies the cryptos.xlsx data refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the edited cells to remain plain text so values such as "1.000"
# or "919.40" are not reinterpreted as numbers and lose their formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.607.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.876.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.84"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2903"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06482"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.99"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07732"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7398"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.876.99"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.93"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.183"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.60"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.599.69"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007476"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.124.52"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.222"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.173"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "165.54"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.192"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.74"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.902"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09895"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.345"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.509"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.247"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.090"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04778"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6938"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01854"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.760"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.227"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.34"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.969"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4161"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8346"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.54"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.336"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.31"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.965"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "919.40"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05670"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.76%  "
